$wb = $excel.ActiveWorkbook
$deac = $wb.Worksheets.Item("deac")
$newSheet = $wb.Worksheets.Add($null, $deac)
$newSheet.Name = "ac"
$newSheet.Range("A1").Value = "work_class"
$newSheet.Range("B1").Value = "ac_cost"
$newSheet.Range("C1").Value = "dmc"
$newSheet.Range("D1").Value = "dollar_basis"

$newSheet.Range("A2").Value = "nohaul"
$newSheet.Range("A2").Style = "Normal"
$newSheet.Range("B2").Formula = "=C2*Markup"
$newSheet.Range("B2").Style = "Output"
$newSheet.Range("C2").Value = 89.445493393168007
$newSheet.Range("C2").Style = "Input"
$newSheet.Range("D2").Value = 2010
$newSheet.Range("D2").Style = "Input"

Write-Host "done"
